$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K holds the 2022 figures, mirroring the layout/formatting of
# the existing 2021 column (J). Copy J's formatting into K first, then
# overwrite the values with the new 2022 data.
$ws.Range("J4:J14").Copy()
$ws.Range("K4:K14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K4").Value = 2022
$ws.Range("K5").Value = 1.6
$ws.Range("K6").Value = 0.4
$ws.Range("K7").Value = 0.9
$ws.Range("K8").Value = 0.6
$ws.Range("K9").Value = 2.1
$ws.Range("K10").Value = 0.6
$ws.Range("K11").Value = 0.9
$ws.Range("K12").Value = 2.3
$ws.Range("K13").Value = 4.3
$ws.Range("K14").Value = 0.3

# Match the recorded selection left by the edit.
$ws.Range("L7").Select()
